# Weekly refresh of the "Zapallo italiano" (Mapocho Venta Directa de Santiago)
# price series. The daily records in rows 2-15 get reshuffled to line up
# with the latest weekly pull: each row's Fecha / Volumen / Precio minimo /
# Precio maximo / Precio promedio ponderado / Unidad de comercializacion /
# Origen / Precio $/Kg / Kg o Unidades move to a (possibly) different row.
# Row 7 is untouched by this refresh.
#
# Columns: D=4 (Fecha), J=10 (Volumen), K=11 (Precio minimo),
#          L=12 (Precio maximo), M=13 (Precio promedio ponderado),
#          N=14 (Unidad de comercializacion), O=15 (Origen),
#          P=16 (Precio $/Kg), Q=17 (Kg o Unidades)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2;  Fecha=44585; Volumen=30; PrecioMin=11000; PrecioMax=11000; PrecioProm=11000; PrecioKg=183 },
    @{ Row=3;  Fecha=44277; Volumen=25 },
    @{ Row=4;  Fecha=44200; Volumen=10; PrecioMin=9000;  PrecioMax=9000;  PrecioProm=9000;  PrecioKg=150 },
    @{ Row=5;  Fecha=45001; Volumen=40; PrecioMin=10000; PrecioMax=10000; PrecioProm=10000; PrecioKg=167 },
    @{ Row=6;  Fecha=44312; Volumen=30; PrecioMax=10000; PrecioProm=10000; PrecioKg=167; Origen="Provincia de Limarí" },
    @{ Row=8;  Fecha=44315; Volumen=25; PrecioMin=10000; PrecioMax=10000; PrecioProm=10000; PrecioKg=167 },
    @{ Row=9;  Fecha=44243; Volumen=80; PrecioMin=10000; PrecioMax=11000; PrecioProm=10375; PrecioKg=173; Origen="Provincia de Quillota" },
    @{ Row=10; Fecha=44284; Volumen=35; PrecioMin=10000; PrecioMax=10000; PrecioProm=10000; PrecioKg=167 },
    @{ Row=11; Fecha=44179; Volumen=15; PrecioMin=7000;  PrecioMax=7000;  PrecioProm=7000;  PrecioKg=117; Unidad="$/caja 60 unidades"; Origen="Provincia de Limarí"; KgUnidades=60 },
    @{ Row=12; Fecha=45030; Volumen=50; PrecioMin=6000;  PrecioMax=6000;  PrecioProm=6000;  PrecioKg=120; Unidad="$/caja 50 unidades"; Origen="Región de Arica y Parinacota"; KgUnidades=50 },
    @{ Row=13; Fecha=44405; Volumen=45; PrecioMin=9000;  PrecioMax=9000;  PrecioProm=9000;  PrecioKg=180; Unidad="$/caja 50 unidades"; Origen="Provincia de Quillota"; KgUnidades=50 },
    @{ Row=14; Fecha=44186; Volumen=15; PrecioMin=7000;  PrecioMax=7000;  PrecioProm=7000;  PrecioKg=117; Unidad="$/caja 60 unidades"; Origen="Provincia de Limarí"; KgUnidades=60 },
    @{ Row=15; Fecha=44291; Volumen=20; PrecioMin=9000;  PrecioMax=9000;  PrecioProm=9000;  PrecioKg=150 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.Fecha

    $ws.Cells.Item($r.Row, 10).Value = $r.Volumen

    if ($r.ContainsKey("PrecioMin")) {
        $ws.Cells.Item($r.Row, 11).Value = $r.PrecioMin
    }
    if ($r.ContainsKey("PrecioMax")) {
        $ws.Cells.Item($r.Row, 12).Value = $r.PrecioMax
    }
    if ($r.ContainsKey("PrecioProm")) {
        $ws.Cells.Item($r.Row, 13).Value = $r.PrecioProm
    }
    if ($r.ContainsKey("Unidad")) {
        $ws.Cells.Item($r.Row, 14).Value = $r.Unidad
    }
    if ($r.ContainsKey("Origen")) {
        $ws.Cells.Item($r.Row, 15).Value = $r.Origen
    }
    if ($r.ContainsKey("PrecioKg")) {
        $ws.Cells.Item($r.Row, 16).Value = $r.PrecioKg
    }
    if ($r.ContainsKey("KgUnidades")) {
        $ws.Cells.Item($r.Row, 17).Value = $r.KgUnidades
    }
}

Write-Output "Weekly row refresh applied."
